# Auto-generated Excel COM-interop script applying the scheduled market-data refresh
# described by the commit diff (per-row currentAveragePrice* / Leve* value updates).
$wb = $excel.ActiveWorkbook

# ALC!row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 9236.5  # H62: 6767.1113 -> 9236.5
$ws.Cells.Item(62, 9).Value = 6000  # I62: 2871.25 -> 6000
$ws.Cells.Item(62, 11).Value = 6000  # K62: 2871.25 -> 6000
$ws.Cells.Item(62, 13).Value = -5376  # M62: -2247.25 -> -5376

# ALC!row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(65, 8).Value = 9236.5  # H65: 6767.1113 -> 9236.5
$ws.Cells.Item(65, 9).Value = 6000  # I65: 2871.25 -> 6000
$ws.Cells.Item(65, 11).Value = 30000  # K65: 14356.25 -> 30000
$ws.Cells.Item(65, 13).Value = -26880  # M65: -11236.25 -> -26880

# ALC!row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 1024.75  # H98: 1218.8 -> 1024.75
$ws.Cells.Item(98, 10).Value = 1999  # J98: 1997 -> 1999
$ws.Cells.Item(98, 12).Value = 1999  # L98: 1997 -> 1999
$ws.Cells.Item(98, 14).Value = -4995  # N98: -4993 -> -4995

# ALC!row 103
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(103, 8).Value = 2198.1538  # H103: 2279.0833 -> 2198.1538
$ws.Cells.Item(103, 9).Value = 2847  # I103: 3171 -> 2847
$ws.Cells.Item(103, 11).Value = 8541  # K103: 9513 -> 8541
$ws.Cells.Item(103, 13).Value = -7955  # M103: -8927 -> -7955

# ALC!row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 1024.75  # H122: 1218.8 -> 1024.75
$ws.Cells.Item(122, 10).Value = 1999  # J122: 1997 -> 1999
$ws.Cells.Item(122, 12).Value = 5997  # L122: 5991 -> 5997
$ws.Cells.Item(122, 14).Value = -10897  # N122: -10891 -> -10897

# ALC!row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 2718.742  # H137: 2768.9312 -> 2718.742
$ws.Cells.Item(137, 9).Value = 2015.6316  # I137: 1961.8889 -> 2015.6316
$ws.Cells.Item(137, 10).Value = 3832  # J137: 4089.5454 -> 3832
$ws.Cells.Item(137, 11).Value = 6046.8948  # K137: 5885.6667 -> 6046.8948
$ws.Cells.Item(137, 12).Value = 11496  # L137: 12268.6362 -> 11496
$ws.Cells.Item(137, 13).Value = -3496.8948  # M137: -3335.6667 -> -3496.8948
$ws.Cells.Item(137, 14).Value = -16596  # N137: -17368.6362 -> -16596

# ALC!row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 2699.5  # H138: 2574.75 -> 2699.5
$ws.Cells.Item(138, 9).Value = 0  # I138: 2433.3333 -> 0
$ws.Cells.Item(138, 10).Value = 2699.5  # J138: 2999 -> 2699.5
$ws.Cells.Item(138, 11).Value = 0  # K138: 7299.999899999999 -> 0
$ws.Cells.Item(138, 12).Value = 8098.5  # L138: 8997 -> 8098.5
$ws.Cells.Item(138, 13).ClearContents()  # M138: -2159.999899999999 -> (removed)
$ws.Cells.Item(138, 14).Value = -18378.5  # N138: -19277 -> -18378.5

# ARM!row 54
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(54, 8).Value = 0  # H54: 10049 -> 0
$ws.Cells.Item(54, 10).Value = 0  # J54: 10049 -> 0
$ws.Cells.Item(54, 12).Value = 0  # L54: 10049 -> 0
$ws.Cells.Item(54, 14).ClearContents()  # N54: -11587 -> (removed)

# ARM!row 58
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(58, 8).Value = 24999.111  # H58: 24998.9 -> 24999.111
$ws.Cells.Item(58, 10).Value = 24999.111  # J58: 24998.9 -> 24999.111
$ws.Cells.Item(58, 12).Value = 24999.111  # L58: 24998.9 -> 24999.111
$ws.Cells.Item(58, 14).Value = -25859.111  # N58: -25858.9 -> -25859.111

# ARM!row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 2222.25  # H74: 2336.2173 -> 2222.25
$ws.Cells.Item(74, 9).Value = 1906.7  # I74: 2028.0526 -> 1906.7
$ws.Cells.Item(74, 11).Value = 1906.7  # K74: 2028.0526 -> 1906.7
$ws.Cells.Item(74, 13).Value = -1032.7  # M74: -1154.0526 -> -1032.7

# ARM!row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 2222.25  # H77: 2336.2173 -> 2222.25
$ws.Cells.Item(77, 9).Value = 1906.7  # I77: 2028.0526 -> 1906.7
$ws.Cells.Item(77, 11).Value = 9533.5  # K77: 10140.263 -> 9533.5
$ws.Cells.Item(77, 13).Value = -5165.5  # M77: -5772.262999999999 -> -5165.5

# ARM!row 131
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(131, 8).Value = 69989  # H131: 69994.5 -> 69989
$ws.Cells.Item(131, 10).Value = 69989  # J131: 69994.5 -> 69989
$ws.Cells.Item(131, 12).Value = 69989  # L131: 69994.5 -> 69989
$ws.Cells.Item(131, 14).Value = -80069  # N131: -80074.5 -> -80069

# BSM!row 80
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 979.4286  # H80: 940.125 -> 979.4286
$ws.Cells.Item(80, 9).Value = 559.3333  # I80: 574.4286 -> 559.3333
$ws.Cells.Item(80, 11).Value = 559.3333  # K80: 574.4286 -> 559.3333
$ws.Cells.Item(80, 13).Value = 438.6667  # M80: 423.5714 -> 438.6667

# BSM!row 83
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(83, 8).Value = 979.4286  # H83: 940.125 -> 979.4286
$ws.Cells.Item(83, 9).Value = 559.3333  # I83: 574.4286 -> 559.3333
$ws.Cells.Item(83, 11).Value = 2796.6665  # K83: 2872.143 -> 2796.6665
$ws.Cells.Item(83, 13).Value = 2195.3335  # M83: 2119.857 -> 2195.3335

# BSM!row 95
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(95, 8).Value = 21999.5  # H95: 25769.834 -> 21999.5
$ws.Cells.Item(95, 10).Value = 21999.5  # J95: 25769.834 -> 21999.5
$ws.Cells.Item(95, 12).Value = 21999.5  # L95: 25769.834 -> 21999.5
$ws.Cells.Item(95, 14).Value = -27491.5  # N95: -31261.834 -> -27491.5

# BSM!row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 1760  # H99: 1497.8 -> 1760
$ws.Cells.Item(99, 9).Value = 1790  # I99: 1426.6666 -> 1790
$ws.Cells.Item(99, 10).Value = 1700  # J99: 1604.5 -> 1700
$ws.Cells.Item(99, 11).Value = 1790  # K99: 1426.6666 -> 1790
$ws.Cells.Item(99, 12).Value = 1700  # L99: 1604.5 -> 1700
$ws.Cells.Item(99, 13).Value = -292  # M99: 71.33339999999998 -> -292
$ws.Cells.Item(99, 14).Value = -4696  # N99: -4600.5 -> -4696

# BSM!row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 2246.25  # H105: 2079.9443 -> 2246.25
$ws.Cells.Item(105, 9).Value = 2331.182  # I105: 2087.8462 -> 2331.182
$ws.Cells.Item(105, 11).Value = 2331.182  # K105: 2087.8462 -> 2331.182
$ws.Cells.Item(105, 13).Value = -584.1819999999998  # M105: -340.8462 -> -584.1819999999998

# CRP!row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 3499.9  # H7: 3387.0645 -> 3499.9
$ws.Cells.Item(7, 9).Value = 5775.6113  # I7: 5471.737 -> 5775.6113
$ws.Cells.Item(7, 11).Value = 5775.6113  # K7: 5471.737 -> 5775.6113
$ws.Cells.Item(7, 13).Value = -5662.6113  # M7: -5358.737 -> -5662.6113

# CRP!row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 3198.25  # H58: 3175.6667 -> 3198.25
$ws.Cells.Item(58, 9).Value = 2448.25  # I58: 2557.6 -> 2448.25
$ws.Cells.Item(58, 11).Value = 2448.25  # K58: 2557.6 -> 2448.25
$ws.Cells.Item(58, 13).Value = -2245.25  # M58: -2354.6 -> -2245.25

# CRP!row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 1604.4286  # H122: 1705.3334 -> 1604.4286
$ws.Cells.Item(122, 10).Value = 2336.3333  # J122: 3005 -> 2336.3333
$ws.Cells.Item(122, 12).Value = 7008.999899999999  # L122: 9015 -> 7008.999899999999
$ws.Cells.Item(122, 14).Value = -11908.9999  # N122: -13915 -> -11908.9999

# CRP!row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 3198.25  # H136: 3175.6667 -> 3198.25
$ws.Cells.Item(136, 9).Value = 2448.25  # I136: 2557.6 -> 2448.25
$ws.Cells.Item(136, 11).Value = 7344.75  # K136: 7672.799999999999 -> 7344.75
$ws.Cells.Item(136, 13).Value = -4794.75  # M136: -5122.799999999999 -> -4794.75

# CUL!row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 98  # H2: 103.5 -> 98
$ws.Cells.Item(2, 10).Value = 30.3  # J2: 32.555557 -> 30.3
$ws.Cells.Item(2, 12).Value = 181.8  # L2: 195.333342 -> 181.8
$ws.Cells.Item(2, 14).Value = -407.8  # N2: -421.333342 -> -407.8

# CUL!row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 151.625  # H12: 161.66667 -> 151.625
$ws.Cells.Item(12, 9).Value = 101.333336  # I12: 121.4 -> 101.333336
$ws.Cells.Item(12, 11).Value = 304.000008  # K12: 364.2 -> 304.000008
$ws.Cells.Item(12, 13).Value = -131.000008  # M12: -191.2 -> -131.000008

# CUL!row 20
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(20, 8).Value = 140  # H20: 0 -> 140
$ws.Cells.Item(20, 10).Value = 140  # J20: 0 -> 140
$ws.Cells.Item(20, 12).Value = 420  # L20: 0 -> 420
$ws.Cells.Item(20, 14).Value = -874  # N20: None -> -874

# CUL!row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 783  # H34: 799.6923 -> 783
$ws.Cells.Item(34, 10).Value = 1074.875  # J34: 1066.5555 -> 1074.875
$ws.Cells.Item(34, 12).Value = 3224.625  # L34: 3199.6665 -> 3224.625
$ws.Cells.Item(34, 14).Value = -3392.625  # N34: -3367.6665 -> -3392.625

# CUL!row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(121, 8).Value = 1556.5  # H121: 1490.6154 -> 1556.5
$ws.Cells.Item(121, 10).Value = 2216.625  # J121: 2048.111 -> 2216.625
$ws.Cells.Item(121, 12).Value = 6649.875  # L121: 6144.333 -> 6649.875
$ws.Cells.Item(121, 14).Value = -9269.875  # N121: -8764.332999999999 -> -9269.875

# CUL!row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 2212.5  # H132: 2149.375 -> 2212.5
$ws.Cells.Item(132, 9).Value = 2040  # I132: 1939 -> 2040
$ws.Cells.Item(132, 11).Value = 18360  # K132: 17451 -> 18360
$ws.Cells.Item(132, 13).Value = -15830  # M132: -14921 -> -15830

# GSM!row 63
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(63, 8).Value = 49995  # H63: 48247 -> 49995
$ws.Cells.Item(63, 9).Value = 0  # I63: 39500 -> 0
$ws.Cells.Item(63, 10).Value = 49995  # J63: 49996.4 -> 49995
$ws.Cells.Item(63, 11).Value = 0  # K63: 39500 -> 0
$ws.Cells.Item(63, 12).Value = 49995  # L63: 49996.4 -> 49995
$ws.Cells.Item(63, 13).ClearContents()  # M63: -38814 -> (removed)
$ws.Cells.Item(63, 14).Value = -51367  # N63: -51368.4 -> -51367

# GSM!row 66
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(66, 8).Value = 49995  # H66: 48247 -> 49995
$ws.Cells.Item(66, 9).Value = 0  # I66: 39500 -> 0
$ws.Cells.Item(66, 10).Value = 49995  # J66: 49996.4 -> 49995
$ws.Cells.Item(66, 11).Value = 0  # K66: 118500 -> 0
$ws.Cells.Item(66, 12).Value = 149985  # L66: 149989.2 -> 149985
$ws.Cells.Item(66, 13).ClearContents()  # M66: -115068 -> (removed)
$ws.Cells.Item(66, 14).Value = -156849  # N66: -156853.2 -> -156849

# GSM!row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1000  # H97: 999.5 -> 1000
$ws.Cells.Item(97, 10).Value = 1000  # J97: 999.5 -> 1000
$ws.Cells.Item(97, 12).Value = 1000  # L97: 999.5 -> 1000
$ws.Cells.Item(97, 14).Value = -1992  # N97: -1991.5 -> -1992

# GSM!row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 2148.5  # H122: 2198 -> 2148.5
$ws.Cells.Item(122, 9).Value = 2148.5  # I122: 2198 -> 2148.5
$ws.Cells.Item(122, 11).Value = 6445.5  # K122: 6594 -> 6445.5
$ws.Cells.Item(122, 13).Value = -3995.5  # M122: -4144 -> -3995.5

# LTW!row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 930  # H22: 954.4 -> 930
$ws.Cells.Item(22, 9).Value = 928.5  # I22: 943 -> 928.5
$ws.Cells.Item(22, 10).Value = 933  # J22: 1000 -> 933
$ws.Cells.Item(22, 11).Value = 928.5  # K22: 943 -> 928.5
$ws.Cells.Item(22, 12).Value = 933  # L22: 1000 -> 933
$ws.Cells.Item(22, 13).Value = -633.5  # M22: -648 -> -633.5
$ws.Cells.Item(22, 14).Value = -1523  # N22: -1590 -> -1523

# LTW!row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 930  # H27: 954.4 -> 930
$ws.Cells.Item(27, 9).Value = 928.5  # I27: 943 -> 928.5
$ws.Cells.Item(27, 10).Value = 933  # J27: 1000 -> 933
$ws.Cells.Item(27, 11).Value = 928.5  # K27: 943 -> 928.5
$ws.Cells.Item(27, 12).Value = 933  # L27: 1000 -> 933
$ws.Cells.Item(27, 13).Value = -821.5  # M27: -836 -> -821.5
$ws.Cells.Item(27, 14).Value = -1147  # N27: -1214 -> -1147

# LTW!row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 1800  # H40: 0 -> 1800
$ws.Cells.Item(40, 9).Value = 1800  # I40: 0 -> 1800
$ws.Cells.Item(40, 11).Value = 1800  # K40: 0 -> 1800
$ws.Cells.Item(40, 13).Value = -1664  # M40: None -> -1664

# LTW!row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 835  # H93: 891.6667 -> 835
$ws.Cells.Item(93, 9).Value = 793.75  # I93: 837.5 -> 793.75
$ws.Cells.Item(93, 11).Value = 793.75  # K93: 837.5 -> 793.75
$ws.Cells.Item(93, 13).Value = 454.25  # M93: 410.5 -> 454.25

# LTW!row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 100  # H122: 0 -> 100
$ws.Cells.Item(122, 9).Value = 100  # I122: 0 -> 100
$ws.Cells.Item(122, 11).Value = 300  # K122: 0 -> 300
$ws.Cells.Item(122, 13).Value = 2150  # M122: None -> 2150

# WVR!row 43
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(43, 8).Value = 0  # H43: 12500 -> 0
$ws.Cells.Item(43, 9).Value = 0  # I43: 12500 -> 0
$ws.Cells.Item(43, 11).Value = 0  # K43: 12500 -> 0
$ws.Cells.Item(43, 13).ClearContents()  # M43: -12351 -> (removed)

# WVR!row 50
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(50, 8).Value = 0  # H50: 14996 -> 0
$ws.Cells.Item(50, 10).Value = 0  # J50: 14996 -> 0
$ws.Cells.Item(50, 12).Value = 0  # L50: 14996 -> 0
$ws.Cells.Item(50, 14).ClearContents()  # N50: -16258 -> (removed)

# WVR!row 131
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(131, 8).Value = 71189.89  # H131: 84905 -> 71189.89
$ws.Cells.Item(131, 12).Value = 71189.89  # L131: 84905 -> 71189.89
$ws.Cells.Item(131, 14).Value = -81269.89  # N131: -94985 -> -81269.89
